$ws = $excel.ActiveWorkbook.ActiveSheet

$values = @{
  2 = "bo"
  3 = "sa"
  4 = "sa"
  5 = "sa"
  6 = "bo"
  7 = "bo"
  8 = "sa"
  9 = "sa"
  10 = "sa"
  11 = "bo"
  12 = "sa"
  13 = "bo"
  14 = "sa"
  15 = "bo"
  16 = "bo"
  17 = "bo"
  18 = "sa"
  19 = "sa"
  20 = "sa"
  21 = "bo"
  22 = "sa"
  23 = "bo"
  24 = "sa"
  25 = "sa"
  26 = "bo"
  27 = "sa"
  28 = "sa"
  29 = "bo"
  30 = "sa"
  31 = "sa"
  32 = "sa"
  33 = "sa"
  34 = "bo"
  35 = "sa"
  36 = "bo"
  37 = "sa"
  38 = "bo"
  39 = "bo"
  40 = "bo"
  41 = "sa"
  42 = "sa"
  43 = "sa"
  44 = "sa"
  45 = "bo"
  46 = "sa"
  47 = "bo"
  48 = "sa"
  49 = "bo"
  50 = "bo"
  51 = "bo"
  52 = "sa"
  53 = "sa"
  54 = "bo"
  55 = "sa"
  56 = "sa"
  57 = "sa"
  58 = "bo"
  59 = "sa"
  60 = "sa"
  61 = "bo"
  62 = "sa"
  63 = "sa"
  64 = "sa"
  65 = "bo"
  66 = "sa"
  67 = "sa"
  68 = "bo"
  69 = "sa"
  70 = "sa"
  71 = "bo"
  72 = "sa"
  73 = "sa"
  74 = "bo"
  75 = "bo"
  76 = "sa"
  77 = "sa"
  78 = "bo"
  79 = "sa"
  80 = "sa"
  81 = "bo"
  82 = "sa"
  83 = "bo"
  84 = "bo"
  85 = "sa"
  86 = "sa"
  87 = "sa"
  88 = "bo"
  89 = "sa"
  90 = "bo"
  91 = "sa"
  92 = "sa"
  93 = "bo"
  94 = "sa"
  95 = "bo"
  96 = "sa"
  97 = "bo"
  98 = "bo"
  99 = "sa"
  100 = "sa"
  101 = "bo"
  102 = "sa"
  103 = "bo"
  104 = "bo"
  105 = "sa"
  106 = "sa"
  107 = "sa"
  108 = "bo"
  109 = "sa"
  110 = "sa"
  111 = "bo"
  112 = "bo"
  113 = "bo"
  114 = "bo"
  115 = "sa"
  116 = "sa"
  117 = "bo"
  118 = "sa"
  119 = "bo"
  120 = "sa"
  121 = "sa"
  122 = "bo"
  123 = "sa"
  124 = "bo"
  125 = "sa"
  126 = "sa"
  127 = "sa"
  128 = "bo"
  129 = "sa"
  130 = "sa"
  131 = "bo"
  132 = "sa"
  133 = "bo"
  134 = "bo"
  135 = "sa"
  136 = "sa"
  137 = "sa"
  138 = "bo"
  139 = "sa"
  140 = "bo"
  141 = "sa"
  142 = "bo"
  143 = "sa"
  144 = "bo"
  145 = "sa"
  146 = "bo"
  147 = "sa"
  148 = "sa"
  149 = "sa"
  150 = "bo"
  151 = "sa"
  152 = "sa"
  153 = "bo"
  154 = "sa"
  155 = "sa"
  156 = "sa"
  157 = "bo"
  158 = "sa"
  159 = "bo"
  160 = "sa"
  161 = "sa"
  162 = "bo"
  163 = "bo"
  164 = "bo"
  165 = "bo"
  166 = "sa"
  167 = "bo"
  168 = "sa"
  169 = "sa"
  170 = "bo"
  171 = "sa"
  172 = "sa"
  173 = "sa"
  174 = "sa"
  175 = "sa"
  176 = "sa"
  177 = "sa"
  178 = "bo"
  179 = "sa"
  180 = "bo"
  181 = "bo"
  182 = "sa"
  183 = "bo"
  184 = "sa"
  185 = "sa"
  186 = "sa"
  187 = "bo"
  188 = "sa"
  189 = "bo"
  190 = "sa"
  191 = "bo"
  192 = "sa"
  193 = "sa"
  194 = "bo"
  195 = "sa"
  196 = "sa"
  197 = "bo"
  198 = "sa"
  199 = "sa"
  200 = "bo"
  201 = "sa"
  202 = "bo"
  203 = "sa"
  204 = "bo"
  205 = "sa"
  206 = "bo"
  207 = "sa"
  208 = "sa"
  209 = "bo"
  210 = "sa"
  211 = "sa"
  212 = "sa"
  213 = "bo"
  214 = "sa"
  215 = "bo"
  216 = "sa"
  217 = "sa"
  218 = "bo"
  219 = "sa"
  220 = "bo"
  221 = "sa"
  222 = "bo"
  223 = "sa"
  224 = "bo"
  225 = "sa"
  226 = "sa"
  227 = "bo"
  228 = "sa"
  229 = "bo"
  230 = "bo"
  231 = "bo"
  232 = "sa"
  233 = "sa"
  234 = "sa"
  235 = "bo"
  236 = "bo"
  237 = "bo"
  238 = "sa"
  239 = "sa"
  240 = "bo"
  241 = "sa"
  242 = "bo"
  243 = "sa"
  244 = "bo"
  245 = "sa"
  246 = "bo"
  247 = "sa"
  248 = "sa"
  249 = "bo"
  250 = "sa"
  251 = "bo"
  252 = "sa"
  253 = "sa"
  254 = "bo"
  255 = "sa"
  256 = "sa"
  257 = "bo"
  258 = "sa"
  259 = "sa"
  260 = "bo"
  261 = "sa"
  262 = "sa"
  263 = "bo"
  264 = "bo"
  265 = "bo"
  266 = "bo"
  267 = "sa"
  268 = "bo"
  269 = "sa"
  270 = "bo"
  271 = "sa"
  272 = "bo"
  273 = "sa"
  274 = "bo"
  275 = "sa"
  276 = "bo"
  277 = "sa"
  278 = "sa"
  279 = "sa"
  280 = "bo"
  281 = "bo"
  282 = "sa"
  283 = "bo"
  284 = "sa"
  285 = "bo"
  286 = "bo"
  287 = "sa"
  288 = "bo"
  289 = "bo"
  290 = "sa"
  291 = "bo"
  292 = "sa"
  293 = "sa"
  294 = "sa"
  295 = "bo"
  296 = "sa"
  297 = "sa"
  298 = "bo"
  299 = "sa"
  300 = "sa"
  301 = "bo"
  302 = "sa"
  303 = "sa"
  304 = "bo"
  305 = "sa"
  306 = "sa"
  307 = "bo"
  308 = "sa"
  309 = "sa"
  310 = "bo"
  311 = "sa"
  312 = "sa"
  313 = "bo"
  314 = "sa"
  315 = "bo"
  316 = "sa"
  317 = "bo"
  318 = "bo"
  319 = "bo"
  320 = "sa"
  321 = "sa"
  322 = "sa"
  323 = "sa"
  324 = "sa"
  325 = "sa"
  326 = "bo"
  327 = "sa"
  328 = "sa"
  329 = "bo"
  330 = "sa"
  331 = "bo"
  332 = "sa"
  333 = "sa"
  334 = "bo"
  335 = "sa"
  336 = "sa"
  337 = "bo"
  338 = "sa"
  339 = "bo"
  340 = "sa"
  341 = "sa"
}

foreach ($row in $values.Keys) {
  $ws.Range("F$row").Value = $values[$row]
}
